$d = $word.ActiveDocument

# --- 1. Resize the columns of the "Fields" table (first table in the doc) ---
# Target gridCol widths (twips): 2131, 935, 3575, 1278  ->  points = twips / 20
$fieldsTable = $d.Tables.Item(1)
$fieldsTable.Columns.Item(1).Width = 106.55
$fieldsTable.Columns.Item(2).Width = 46.75
$fieldsTable.Columns.Item(3).Width = 178.75
$fieldsTable.Columns.Item(4).Width = 63.9

# --- 2. "When the comment was added to the application." -> "...SPPA." ---
$d.Content.Find.Execute("added to the application.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "added to the SPPA.", 2)

# --- 3. "can be seen by all parties involved in the application." -> "...SPPA." ---
$d.Content.Find.Execute("involved in the application.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "involved in the SPPA.", 2)
